# Swap the "category-code" (column F) and "category-name" (column G) values
# on every row of the sheet, including the header row, to correct the
# column ordering bug described in the commit (codeforIATI SectorGroup
# codelist regeneration).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row/column so this works regardless of how many
# data rows are present.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$firstRowIndex = $usedRange.Row
$lastRowIndex = $firstRowIndex + $lastRow - 1

for ($r = $firstRowIndex; $r -le $lastRowIndex; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    # Read through .Text so numeric-looking strings (e.g. "111") come back
    # as their original text rather than being coerced to a Double (which
    # .Value/.Value2 would do).
    $fVal = $fCell.Text
    $gVal = $gCell.Text

    # Re-assign with a leading apostrophe to force text storage, preserving
    # the original (string) cell type even for numeric-looking codes (e.g.
    # "111") which would otherwise be auto-coerced to a number.
    $fCell.Value = "'" + $gVal
    $gCell.Value = "'" + $fVal

    # Restore the default "Normal" style so the apostrophe/quote-prefix
    # bookkeeping doesn't leave a visible style change on the cell itself.
    $fCell.Style = "Normal"
    $gCell.Style = "Normal"
}
